# Add a new "UK" worksheet, cloned from the "Poland" worksheet, and
# populate it with the UK-specific test data (mirrors the commit
# "Added Test Data for UK Market").

$wb = $excel.ActiveWorkbook

# The "Poland" sheet is the template for every per-market sheet in this
# workbook (same layout/styles/merged cells), so copy it and place the
# copy at the very end of the tab strip.
$template = $wb.Worksheets.Item("Poland")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)

# The copy becomes the last sheet; rename it and fill in the UK values.
$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# Set the user-story code cell first, then the market name cell, so the
# new shared-string entries are appended in the same order Excel used.
$uk.Range("B4").Value = "NGC-2741/T3345"
$uk.Range("B2").Value = "UK Market"

# Match the author's final selection on the new sheet.
$uk.Range("B4").Select()
